$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.974.78"
$ws.Range("E2").Value = "  +5.20%  "
$ws.Range("D3").Value = "1.881.06"
$ws.Range("E3").Value = "  +4.19%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'282.70"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "'0.5269"
$ws.Range("E7").Value = "  +3.70%  "
$ws.Range("D8").Value = "'0.3546"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("D9").Value = "'45.27"
$ws.Range("E9").Value = "  +3.72%  "
$ws.Range("D10").Value = "'0.07085"
$ws.Range("E10").Value = "  +6.68%  "
$ws.Range("D11").Value = "'20.40"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").Value = "'0.8203"
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("D13").Value = "'0.07830"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("D14").Value = "1.889.03"
$ws.Range("E14").Value = "  +4.61%  "
$ws.Range("D15").Value = "'5.238"
$ws.Range("E15").Value = "  +3.35%  "
$ws.Range("E16").Value = "  +3.73%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "'14.58"
$ws.Range("E18").Value = "  +4.72%  "
$ws.Range("D19").Value = "'0.000008176"
$ws.Range("E19").Value = "  +2.92%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "27.045.14"
$ws.Range("E21").Value = "  +5.25%  "
$ws.Range("D22").Value = "'4.800"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").Value = "'6.274"
$ws.Range("E24").Value = "  +3.99%  "
$ws.Range("D25").Value = "'2.413"
$ws.Range("E25").Value = "  +14.36%  "
$ws.Range("D26").Value = "'146.97"
$ws.Range("E26").Value = "  +3.12%  "
$ws.Range("D27").Value = "'17.68"
$ws.Range("E27").Value = "  +4.68%  "
$ws.Range("D28").Value = "'1.667"
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("D29").Value = "'114.10"
$ws.Range("E29").Value = "  +5.33%  "
$ws.Range("D30").Value = "'4.415"
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("D31").Value = "'4.401"
$ws.Range("E31").Value = "  +4.33%  "
$ws.Range("D32").Value = "'0.08891"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").Value = "'0.04926"
$ws.Range("E33").Value = "  +2.83%  "
$ws.Range("D34").Value = "'1.181"
$ws.Range("E34").Value = "  +5.00%  "
$ws.Range("D35").Value = "'0.7498"
$ws.Range("E35").Value = "  +3.73%  "
$ws.Range("D36").Value = "'2.896"
$ws.Range("E36").Value = "  +1.44%  "
$ws.Range("D37").Value = "'3.298"
$ws.Range("E37").Value = "  +8.68%  "
$ws.Range("D38").Value = "'2.409"
$ws.Range("E38").Value = "  +5.00%  "
$ws.Range("D39").Value = "'0.5329"
$ws.Range("E39").Value = "  +3.08%  "
$ws.Range("D40").Value = "'0.01897"
$ws.Range("E40").Value = "  +2.14%  "
$ws.Range("D41").Value = "'0.9837"
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("D42").Value = "'117.18"
$ws.Range("E42").Value = "  +2.09%  "
$ws.Range("D43").Value = "'6.325"
$ws.Range("E43").Value = "  +2.32%  "
$ws.Range("D44").Value = "'8.215"
$ws.Range("E44").Value = "  +2.36%  "
$ws.Range("D45").Value = "'0.4644"
$ws.Range("E45").Value = "  +1.42%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").Value = "'9.488"
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D49").Value = "'36.89"
$ws.Range("E49").Value = "  +2.97%  "
$ws.Range("D50").Value = "'1.531"
$ws.Range("E50").Value = "  +2.72%  "
$ws.Range("D51").Value = "'0.05949"
$ws.Range("E51").Value = "  +2.33%  "
